$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text, matching the source format
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.070.68"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.857.88"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "237.68"
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("E6").Value = "  +1.59%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "42.36"
$ws.Range("E8").Value = "  +8.96%  "
$ws.Range("D9").Value = "0.329"
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").Value = "0.0697"
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.862.30"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").Value = "11.40"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "0.677"
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("D16").Value = "4.70"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("D17").Value = "35.030.98"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "70.28"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").Value = "240.94"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "12.16"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").Value = "171.59"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "1.95"
$ws.Range("E26").Value = "  +31.24%  "
$ws.Range("D27").Value = "7.91"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("D28").Value = "17.71"
$ws.Range("E28").Value = "  +3.54%  "
$ws.Range("D29").Value = "0.125"
$ws.Range("E29").Value = "  +3.33%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.0558"
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("B31").Value = "BinanceUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D31").Value = "1.01"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "4.00"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "4.01"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Value = "2.03"
$ws.Range("E34").Value = "  +13.86%  "
$ws.Range("E35").Value = "  +22.38%  "
$ws.Range("D36").Value = "0.787"
$ws.Range("E36").Value = "  +14.24%  "
$ws.Range("D37").Value = "1.29"
$ws.Range("E37").Value = "  +5.00%  "
$ws.Range("E38").Value = "  +12.79%  "
$ws.Range("D39").Value = "91.74"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  +6.72%  "
$ws.Range("D41").Value = "1.353.57"
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("D42").Value = "14.87"
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  +6.59%  "
$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "12.74"
$ws.Range("E44").Value = "  +57.53%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").Value = "0.0544"
$ws.Range("E47").Value = "  +6.06%  "
$ws.Range("D48").Value = "6.37"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("E50").Value = "  +3.25%  "
$ws.Range("D51").Value = "3.42"
